$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('14', 'Dakka 1 gol 1 hemen çöktü güzel okulumuzun güzel ders seçme sayfası yahu sadece 4. Sınıflar seçiyor ama yine de çöküyor altı üstü kaç bin kişi var da düzgün bi site yapamıyorlar'),
    @('15', 'Ders programı hala açıklanmadı. Son sınıfım bir sürü seçmeli dersim var hangisini alacağımı bilmiyorum ve yarın ders seçimi var. Okulu arıyorum yardımcı olan yok. Şaşırtmadın yine Işık Üniversitesi.'),
    @('16', 'maslak kampüsünün şileye taşınmasıyla birlikte yurt kapasitesinin çok üstünde bir başvuru olmuş. bu da demek oluyor ki başvuru yapan herkes kabul edilemeyecek ve açıkta kalacak. ayrıca yurt sonuçları okulun açılmasından 1 hafa önce belli olacak yani hala belli değil. diyelim ki yurt başvurumuz kabul edilmedi sokakta mı yatmamız gerekiyor?! madem okulda yeterince yurt kapasitesi yok o zaman kampüs taşımanın acelesi mi vardı? hadi taşıdın diyelim bu sorunlar niye gözardı ediliyo? bir alternatif düşünülmemiş ve öğrencilerin de bir alternatif plan yapmasına fırsat bile tanınmıyor. 1 hafta içerisinde insanlar nasıl başka bir yer bulsun? o kadar heyet kurul bilmem ne bu kadar karar alıyor da birinin de aklına bu problemler gelmiyo mu? hadi diyelim geldi ve gözardı ettiler, e bari düzgün bir açıklama yapın, ihtimallerden bahsedin, en önemlisi erken davranın..'),
    @('17', 'Mali işlere yazıyorum öğrenci işlerine aktarıyor oraya yazıyorum başka bi yere aktarıyor. Telefonlar kesinlikle açılmıyor . Sorunlar kesinlikle çözülmüyor. Mahalle bakkalım daha iş bitirici bu okul çalışanlarından. Yazık'),
    @('18', 'Yani bu nasıl mantık ben anlamış değilim. Maslaktan Şileye gelen öğrenciye %50 indirim yapıldığı zaman yurtlarda, zaten halihazırda şilede okuyan, yurtlarda kalmak durumunda olan öğrencilere haksızlık olmayacak mı? Bizler niye her lafa boyun eğiyoruz yetişkin insanlar olarak niye mantıksız kararlara itiraz etmiyoruz? Ya da en basitinden öğrenci konseyi (!) böyle rezil bi kararı paylaşabiliyor? Tamam maslaktakilerin de kampüs hayatı yaşayacak olması çok iyi onlar adına seviniyorum ama kimsenin hakkı gasp olmadan olmalı bu iş. Umarım konsey itiraf sayfalarını okuyordur ve yanlışlarını görüyordur.'),
    @('19', 'Hocalarımız neden Final notlarını açıklamamakta diretiyorlar? Neden biz finalden kaç aldığımızı öğrenemiyoruz? Mail atıyorum Hoca harf notlarını girdim finalleri açıklamıyorum diyor yani amaç ne ?'),
    @('20', 'Ağustos''ta ödeyip ödeme dekontunu mail atıp teyidini aldığım yurt borcumu sistemden silmedikleri ve 2 gündür ne bir mailime ne de bir telefonuma cevap vermedikleri için ders seçemedim. Teşekkürler ışık üniversitesi, harika bir üniversite deneyimi yaşatıyorsun gerçekten.'),
    @('21', 'Işık üniversitesi yine öğrencisiyle dalga geçmeye devam ediyor. İki fakülteyi şileye taşı ama iki tane servis ver geç. Öğrencisini hiçe sayan, sadece müşteri olarak gören bir kurum.'),
    @('22', 'Okulun acılmasına aşağı yukarı 1 ay var ve hala online mı yüz yüze mi belli değil. Bu insanlar ne zaman ev/yurt bakacak? Neden acıklama yapmıyorlar? Okul online olursa hem çalışıp hem okuyacak arkadaşlar var, yüz yüze olursa şehir dışından gelecek insanlar var... Neden yani her şeyi bu kadar yokuşa sürmek neden? Okula yeni gelen arkadaşlar; cehenneme hoşgeldiniz.'),
    @('23', 'Demeyim demeyim dedim ama okulda Öğrenci Konseyi yok arkadaşlar. Başkanı kim, ne zaman seçildi, kim seçti bilen var mı? Devir teslim düzeninin olduğu, öğrencilerin seçmediği konseyin (çalışsalar bile) öğrencilere faydası ne kadar olabilir? Büyük hayaller, tarihi, güzel kampüsü hatta sırf adını Atatürk verdiği için gururlandığım, "Işıklı" olmanın ayrıcalık olduğunu düşündüğüm bir başlangıçtan utanarak mezun olma kısmına getirdiler beni. Umarım Türkiye gibi hak ettiği değişimi her alanda görür ve eskiye döner. Kalanlara sabır diliyorum.'),
    @('24', 'Boşuna çene yoruyorsunuz. Burası bir okul değil ticarethane. Okul, pandemi dolayısıyla 2 dönem online oldu, bir ton masraftan kurtuldular ama yine de öğrencileri düşünüp indirim yapmak yerine fiyatları arttırdılar. Sonradan kdv''yi az aldılar. Hocalar desen not açıklamaz, kafasına göre not verir. Shuttle''lar desen git-gel 50 lira. Şimdi de giderleri en aza indirip herkesi bir kampüste toplayıp daha çok kazanacaklar. Bu kararı da şile belediyesi ile karşılıklı çıkar ilişkisi güderek almadılarsa ben de hiçbir şey bilmiyorum. Yeni tercih etmeyi düşünüp de burayı okuyan varsa GELMEYİN kardeşim bu okula. Batsınlar. Ya kaydınızı alın ya da mezun olup kurtulun.'),
    @('25', 'Admin gizli olsun ya ben bu ekonomi hocalarının sınav kurallarını anlamaya çalışıyorum acaba kendileri derslerde slayt okumaktan başka bir iş yapmayıp konuyu adam akıllı öğretememelerinin faturasını bize katı kuralları koyarak mı kesiyorlar anlamadım ki kalem kullanma anladık ulan ses duyduğunda disipline göndermek veya puan kırmak nedir nefes alma diyorlar resmen insana.'),
    @('26', 'Hocalar olarak yapabileceğimiz bir şey yok tavrı takınan politik doğrucu hocalarımıza şunu söylemek istiyorum; en azından bu kararı desteklemediğinizi söyleyin, öğrencilerin mağduriyetinden bahsedin. Öğrencinin mi, maaş çeklerinizin mi dostusunuz oturup bir düşünün.'),
    @('27', 'okulu arayınca mail atmam gerektiğini söylüyolar fakat attığım maile günlerdir cevap vermiyolar.Bu öğrenci işlerine ulaşmak ne kadar zormuş. Başka ne şekil ulaşabilirim ? illa okula mı gitmem lazım ulaşabilmek için.')
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = [int]$data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}

$ws.Range("G33").Select() | Out-Null
